$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) "Fix last image in transfer sample": on the last slide, the marble
#    value was mistakenly shown as 110 in three places; it should be 100
#    (matching the other occurrences already on the slide).
# ---------------------------------------------------------------------------
$lastSlide = $p.Slides.Item($p.Slides.Count)

for ($i = 1; $i -le $lastSlide.Shapes.Count; $i++) {
    $shp = $lastSlide.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }

    $tr = $shp.TextFrame.TextRange
    $full = $tr.Text

    # "hash(110)" -> "hash(100)"
    while ($full.IndexOf("hash(110)") -ge 0) {
        $idx = $full.IndexOf("hash(110)")
        $sub = $tr.Characters($idx + 1, 9)
        $sub.Text = "hash(100)"
        $full = $tr.Text
    }

    # " 110" -> " 100" (the standalone value, not part of hash(110))
    while ($full.IndexOf(" 110") -ge 0) {
        $idx = $full.IndexOf(" 110")
        $sub = $tr.Characters($idx + 1, 4)
        $sub.Text = " 100"
        $full = $tr.Text
    }
}

# ---------------------------------------------------------------------------
# 2) The auto-updating "date" placeholder re-cached to the next day
#    (5/20/20 -> 5/21/20) on every layout + the slide master.
# ---------------------------------------------------------------------------
$master = $p.SlideMaster

function Fix-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if (-not $shp.HasTextFrame) { continue }
        if ($shp.TextFrame.TextRange.Text -eq "5/20/20") {
            $shp.TextFrame.TextRange.Text = "5/21/20"
        }
    }
}

Fix-DatePlaceholder $master

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    Fix-DatePlaceholder $master.CustomLayouts.Item($L)
}
